# Split Lot AJ02, AJ03 into separate records in the bead catalog.
#
# The sheet currently has a single calibration-bead block (rows 77-95,
# fluorophores BV421..BUV737) labelled "Lot AJ02" in column A of its first
# row (A77). This adds a second, identical block for "Lot AJ03" directly
# below the existing data (new rows 96-114), duplicating the B:K content
# (fluorophore name, laser, filter, MEF-unit name, and the six bead-peak
# values) and labelling its first row with the new lot name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the B:K contents (values + number formats) of the existing
# "Lot AJ02" block (rows 77-95) into the new block (rows 96-114). Column A
# is intentionally excluded here so only row 96 ends up with a value in
# column A (set explicitly below) -- rows 97-114 stay blank in column A,
# matching the source block where only row 77 carries the lot label.
$src = $ws.Range("B77:K95")
$dest = $ws.Range("B96:K114")
$src.Copy($dest)

# Label the new block with the new lot number.
$ws.Range("A96").Value = "Lot AJ03"

# Leave the new label cell selected/active, matching the natural end state
# of a user who just typed the new lot name.
[void]$ws.Range("A96").Select()
